$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "test3Arguments" example block (rows 11:12) into a new
# "test4Arguments" block (rows 19:20), reusing the same look & feel
# (fills, borders, blank helper cells, merged label cell A19:A20).
$ws.Range("A11:G12").Copy($ws.Range("A19:G20"))

# Replace the copied text with the new VarArgs test case content.
# Order matches how the strings were appended to the shared-string table.
$ws.Range("C19").Value = "Method int test4Arguments()"
$ws.Range("G19").Value = "Method int test4(Object[] args)"
$ws.Range("G20").Value = "return args.length;"
$ws.Range("C20").Value = "return test4((Integer) 5, ""abc"", (Integer) 10, (Integer) 12, (Double) 14);"

$ws.Range("C23").Select() | Out-Null
